# Update ESTADO column (B) from "PENDIENTE" to "VOLADA" for the specified rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$rows = @(158,159,160,161,162,163,164,165,166,167,168,170,188,189,197,198,199,200,201)

foreach ($r in $rows) {
    $ws.Range("B$r").Value = "VOLADA"
}
